# Restructure the "Import Data Barang" template table on Sheet1.
#
# The original header row had separate "Kredit" and "Cash" pricing columns
# (Kredit Dus, Kredit Pack, Kredit Pcs, Cash Dus, Cash Pack, Cash Pcs) plus
# Diskon and the Min Qty note. The new structure drops the "Kredit" columns
# entirely and also drops "Cash Pcs", leaving:
#   Nama Suplier | Kode Barang | Nama Barang | Cash Dus | Cash Pack | Diskon | Min Qty (...)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Kredit Dus", "Kredit Pack", "Kredit Pcs" columns (D:F).
# This shifts the old Cash Dus/Cash Pack/Cash Pcs/Diskon/Min Qty columns
# (G:K) left to D:H.
$ws.Range("D1:F4").EntireColumn.Delete()

# Remove the "Cash Pcs" column, which is now column F after the shift
# above. This shifts Diskon/Min Qty (G:H) left to F:G.
$ws.Range("F1:F4").EntireColumn.Delete()

# Match the saved selection state of the edited workbook (last active
# cell ends up on the new last header cell, G4).
$ws.Range("G4").Select() | Out-Null
